$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate the weekly blog-post references in row 8:
#   I8 previously showed "ser: 119" -> now shows "ser: 120"
#   E8 previously showed "ser: 120" -> now shows "ser: 121"
#   C8 previously showed "ser: 121" -> now shows "ser: 122" (new post)
# D8 (the meetup card) keeps its own text unchanged.
$ws.Range("I8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 120"
$ws.Range("E8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 121"
$ws.Range("C8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 122"

# Update the active selection to I8.
$ws.Range("I8").Select() | Out-Null
